# Fixed a problem with export from Scholarship Universe.
# The ApplicantRanking-amount column (B) previously pointed at shared-string
# text values ("30000", "15000", "6000", "5000", "1200"); it should instead
# hold literal numbers. At the same time the per-scholarship
# applicant-number / ranking pairs (columns C/D) need to be re-written with
# their corrected (re-shuffled) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 2..62, columns A (ScholarshipName id), B (TotalAmountAvailable),
# C (applicantnumber), D (ApplicantRanking) - final corrected values.
$data = @(
    @(15,30000,20,3),
    @(15,30000,9,8),
    @(15,30000,49,10),
    @(15,30000,3,11),
    @(15,30000,7,17),
    @(15,30000,36,19),
    @(15,30000,31,21),
    @(15,30000,25,23),
    @(15,30000,27,25),
    @(15,30000,42,31),
    @(15,30000,6,32),
    @(15,30000,18,40),
    @(15,30000,34,44),
    @(2,15000,19,2),
    @(2,15000,21,18),
    @(2,15000,46,22),
    @(1,6000,14,6),
    @(1,6000,32,15),
    @(1,6000,36,19),
    @(1,6000,31,21),
    @(1,6000,27,25),
    @(1,6000,42,31),
    @(1,6000,17,33),
    @(1,6000,47,42),
    @(6,6000,14,6),
    @(6,6000,8,28),
    @(8,5000,19,2),
    @(8,5000,14,6),
    @(8,5000,17,33),
    @(12,5000,8,28),
    @(12,5000,1,34),
    @(12,5000,2,36),
    @(12,5000,39,37),
    @(13,5000,19,2),
    @(13,5000,50,16),
    @(14,5000,19,2),
    @(14,5000,14,6),
    @(14,5000,24,30),
    @(14,5000,35,45),
    @(9,1200,19,2),
    @(9,1200,15,5),
    @(9,1200,14,6),
    @(9,1200,44,9),
    @(9,1200,38,12),
    @(9,1200,33,13),
    @(9,1200,32,15),
    @(9,1200,36,19),
    @(9,1200,31,21),
    @(9,1200,25,23),
    @(9,1200,51,24),
    @(9,1200,27,25),
    @(9,1200,23,26),
    @(9,1200,43,27),
    @(9,1200,8,28),
    @(9,1200,42,31),
    @(9,1200,17,33),
    @(9,1200,39,37),
    @(9,1200,29,38),
    @(9,1200,41,39),
    @(9,1200,28,43),
    @(9,1200,34,44)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Reflect the selection stored in the sheet view (A1:D62 selected).
$ws.Range("A1:D62").Select()
